$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "251×7=1757"
$t.Cell(1, 2).Range.Text = "532×7=3724"
$t.Cell(1, 3).Range.Text = "966×3=2898"
$t.Cell(1, 4).Range.Text = "768×9=6912"
$t.Cell(1, 5).Range.Text = "692×7=4844"
$t.Cell(5, 1).Range.Text = "398×9=3582"
$t.Cell(5, 2).Range.Text = "554×4=2216"
$t.Cell(5, 3).Range.Text = "514×6=3084"
$t.Cell(5, 4).Range.Text = "868×3=2604"
$t.Cell(5, 5).Range.Text = "820×5=4100"
$t.Cell(10, 1).Range.Text = "553×4=2212"
$t.Cell(10, 2).Range.Text = "412×3=1236"
$t.Cell(10, 3).Range.Text = "215×8=1720"
$t.Cell(10, 4).Range.Text = "529×7=3703"
$t.Cell(10, 5).Range.Text = "901×3=2703"
$t.Cell(15, 1).Range.Text = "298×5=1490"
$t.Cell(15, 2).Range.Text = "963×5=4815"
$t.Cell(15, 3).Range.Text = "370×7=2590"
$t.Cell(15, 4).Range.Text = "836×4=3344"
$t.Cell(15, 5).Range.Text = "437×8=3496"
$t.Cell(20, 1).Range.Text = "914×9=8226"
$t.Cell(20, 2).Range.Text = "166×5=830"
$t.Cell(20, 3).Range.Text = "184×8=1472"
$t.Cell(20, 4).Range.Text = "228×5=1140"
$t.Cell(20, 5).Range.Text = "892×9=8028"
